$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'300.61"
$ws.Range("E2").Value = "'-4.66%"
$ws.Range("D3").Value = "'35.35"
$ws.Range("E3").Value = "'-0.87%"
$ws.Range("E4").Value = "'-0.99%"
$ws.Range("D5").Value = "'0.07979"
$ws.Range("E5").Value = "'-2.32%"
$ws.Range("D6").Value = "'1.928"
$ws.Range("E6").Value = "'-7.77%"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'7.790"
$ws.Range("E7").Value = "'-1.89%"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "'4.051"
$ws.Range("E8").Value = "'-2.09%"
$ws.Range("D9").Value = "'0.9236"
$ws.Range("E9").Value = "'-1.16%"
$ws.Range("D10").Value = "'0.1493"
$ws.Range("E10").Value = "'43.62%"
$ws.Range("D11").Value = "'0.1893"
$ws.Range("E11").Value = "'-1.42%"
$ws.Range("D12").Value = "'0.09193"
$ws.Range("E12").Value = "'1.38%"
$ws.Range("D13").Value = "'0.03412"
$ws.Range("E13").Value = "'-5.09%"
$ws.Range("D14").Value = "'0.09856"
$ws.Range("E14").Value = "'-0.50%"
$ws.Range("D15").Value = "'0.001390"
$ws.Range("E15").Value = "'-2.80%"
$ws.Range("D16").Value = "'0.005744"
$ws.Range("E16").Value = "'-1.01%"
$ws.Range("D17").Value = "'3.512"
$ws.Range("E17").Value = "'1.21%"
$ws.Range("D18").Value = "'2.991"
$ws.Range("E18").Value = "'0.16%"
$ws.Range("E19").Value = "'-1.61%"
$ws.Range("D20").Value = "'0.1293"
$ws.Range("E20").Value = "'-1.40%"
$ws.Range("E21").Value = "'-1.06%"
$ws.Range("D22").Value = "'0.2403"
$ws.Range("E22").Value = "'8.56%"
$ws.Range("D23").Value = "'0.04462"
$ws.Range("E23").Value = "'-2.03%"
$ws.Range("D24").Value = "'0.001218"
$ws.Range("E24").Value = "'-1.93%"
$ws.Range("D25").Value = "'0.004779"
$ws.Range("E25").Value = "'-0.46%"
$ws.Range("D26").Value = "'0.0001233"
$ws.Range("E26").Value = "'-1.49%"
$ws.Range("D27").Value = "'0.0003005"
$ws.Range("E27").Value = "'-33.31%"
$ws.Range("D39").Value = "'0.01899"
$ws.Range("E39").Value = "'-4.02%"
$ws.Range("D40").Value = "'0.04722"
$ws.Range("E40").Value = "'-4.03%"
$ws.Range("D41").Value = "'0.007377"
$ws.Range("E41").Value = "'-2.93%"
$ws.Range("D42").Value = "'0.009723"
$ws.Range("E42").Value = "'23.44%"
$ws.Range("D43").Value = "'0.1329"
$ws.Range("E43").Value = "'-4.03%"
$ws.Range("D44").Value = "'0.002116"
$ws.Range("E44").Value = "'0.87%"
$ws.Range("D45").Value = "'0.009327"
$ws.Range("E45").Value = "'-20.54%"
$ws.Range("D46").Value = "'0.00006263"
$ws.Range("E46").Value = "'-6.90%"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.12%"
$ws.Range("D48").Value = "'65.00"
$ws.Range("E48").Value = "'-58.57%"
$ws.Range("D49").Value = "'0.001662"
$ws.Range("E49").Value = "'-2.40%"
$ws.Range("D50").Value = "'0.00002106"
$ws.Range("E50").Value = "'0.12%"
$ws.Range("D51").Value = "'0.0002005"
$ws.Range("E51").Value = "'0.12%"
